$wb = $excel.ActiveWorkbook

# Parameterize Browser election in runtime:
# Flip the "runmode" flag for the last AddCustomerTest data row ("Bassel Safwat")
# from "Y" to "N" so it is excluded from the run.
$ws = $wb.Worksheets.Item("AddCustomerTest")
$ws.Range("E5").Value = "N"
